## Adds three new booking records (rows 151-153) to the "Modelling" sheet,
## mirroring the formatting/formulas already used for the existing data
## rows, and leaves the selection on AC147 (matching the saved view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Seed rows 151:153 with the formatting used by the last existing data
#    row (150) so borders / fonts / alignment / number formats all line up
#    with the rest of the table before we overwrite the actual values.
# ---------------------------------------------------------------------------
$ws.Range("A150:AP150").Copy()
$ws.Range("A151:AP153").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Columns F / M / N hold real dates on these three rows (row 150 used plain
# text there), so pull the numeric-date style from an existing date cell.
$ws.Range("F7").Copy()
$ws.Range("F151:F153").PasteSpecial(-4122)
$ws.Range("M152:M153").PasteSpecial(-4122)
$ws.Range("N152:N153").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column G holds a real time-of-day value (instead of the text that the
# rest of the sheet uses), so start from the date style too and then swap
# the number format to a time format - this creates the new "h:mm:ss"
# cell style used by the workbook.
$ws.Range("G151:G153").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("G151:G153").NumberFormat = "h:mm:ss"

# ---------------------------------------------------------------------------
# 2. Row 151 values
# ---------------------------------------------------------------------------
$ws.Range("A151").Value = 8239
$ws.Range("B151").Value = 11
$ws.Range("C151").Value = "Business"
$ws.Range("D151").Value = "International"
$ws.Range("E151").Value = "Hotel"
$ws.Range("F151").Value = 43816
$ws.Range("G151").Value = 0.59508101851851858
$ws.Range("H151").Value = "Qatar"
$ws.Range("I151").Value = "Doha"
$ws.Range("J151").Value = "Victoria Hotel"
$ws.Range("K151").Value = 4
$ws.Range("L151").Value = "Confimed"
$ws.Range("M151").Value = "19-12-2018"
$ws.Range("N151").Value = "20-12-2018"
$ws.Range("O151").Value = 1
$ws.Range("P151").Value = 1
$ws.Range("Q151").Value = 1
$ws.Range("R151").Value = 0
$ws.Range("S151").Value = 0
$ws.Range("T151").Value = 0
$ws.Range("U151").Value = "Room Only"
$ws.Range("V151").Value = "Standard Room"
$ws.Range("W151").Value = "Single"
$ws.Range("X151").Value = "Refundable"
$ws.Range("Y151").Value = "Early Arrival"
$ws.Range("Z151").Value = "No Change"
$ws.Range("AA151").Formula = "=AD151-(AD151*15%)"
$ws.Range("AB151").Value = 0
$ws.Range("AC151").Value = 0
$ws.Range("AD151").Value = 15550
$ws.Range("AE151").Value = 15550
$ws.Range("AF151").Formula = "=AG151"
$ws.Range("AG151").Formula = "=AD151-AA151"
$ws.Range("AH151").Value = "Y"
$ws.Range("AI151").Formula = "=IF(O151<=3,""ShortStay"",""LongStay"")"
$ws.Range("AJ151").Value = "Business - Single"
$ws.Range("AK151").Value = 4
$ws.Range("AL151").Value = "Indian"
$ws.Range("AM151").Value = 25
$ws.Range("AN151").Value = "Doha"
$ws.Range("AO151").Value = "Qatar"
$ws.Range("AP151").Value = "Mobile"

# ---------------------------------------------------------------------------
# 3. Row 152 values
# ---------------------------------------------------------------------------
$ws.Range("A152").Value = 8240
$ws.Range("B152").Value = 43
$ws.Range("C152").Value = "Leisure"
$ws.Range("D152").Value = "International"
$ws.Range("E152").Value = "Hotel"
$ws.Range("F152").Value = 43823
$ws.Range("G152").Value = 0.73961805555555549
$ws.Range("H152").Value = "France"
$ws.Range("I152").Value = "Noisy-le-Grand"
$ws.Range("J152").Value = "Novotel Marne La Vallee Noisy Le Grand"
$ws.Range("K152").Value = 4
$ws.Range("L152").Value = "Cancelled"
$ws.Range("M152").Value = 43826
$ws.Range("N152").Value = 43826
$ws.Range("O152").Value = 2
$ws.Range("P152").Value = 1
$ws.Range("Q152").Value = 1
$ws.Range("R152").Value = 0
$ws.Range("S152").Value = 0
$ws.Range("T152").Value = 0
$ws.Range("U152").Value = "BED AND BREAKFAST"
$ws.Range("V152").Value = "Deluxe Room"
$ws.Range("W152").Value = "Single"
$ws.Range("X152").Value = "Refundable"
$ws.Range("Y152").Value = "Early Arrival"
$ws.Range("Z152").Value = "No Change"
$ws.Range("AA152").Formula = "=AD152-(AD152*25%)"
$ws.Range("AB152").Value = 0
$ws.Range("AC152").Value = 191.76
$ws.Range("AD152").Value = 2000
$ws.Range("AE152").Value = 1689.98
$ws.Range("AF152").Formula = "=AG152"
$ws.Range("AG152").Formula = "=AD152-AA152"
$ws.Range("AH152").Value = "Y"
$ws.Range("AI152").Formula = "=IF(O152<=3,""ShortStay"",""LongStay"")"
$ws.Range("AJ152").Value = "Leisure - Single"
$ws.Range("AK152").Value = 4.5
$ws.Range("AL152").Value = "Indian"
$ws.Range("AM152").Value = 25
$ws.Range("AN152").Value = "Vancouver"
$ws.Range("AO152").Value = "Canada"
$ws.Range("AP152").Value = "Mobile"

# ---------------------------------------------------------------------------
# 4. Row 153 values
# ---------------------------------------------------------------------------
$ws.Range("A153").Value = 8241
$ws.Range("B153").Value = 33
$ws.Range("C153").Value = "Leisure"
$ws.Range("D153").Value = "International"
$ws.Range("E153").Value = "Hotel"
$ws.Range("F153").Value = 43824
$ws.Range("G153").Value = 0.27656249999999999
$ws.Range("H153").Value = "Indonesia"
$ws.Range("I153").Value = "Medan"
$ws.Range("J153").Value = "JW Marriott Hotel Medan"
$ws.Range("K153").Value = 5
$ws.Range("L153").Value = "Cancelled"
$ws.Range("M153").Value = 43826
$ws.Range("N153").Value = 43824
$ws.Range("O153").Value = 2
$ws.Range("P153").Value = 1
$ws.Range("Q153").Value = 1
$ws.Range("R153").Value = 0
$ws.Range("S153").Value = 0
$ws.Range("T153").Value = 0
$ws.Range("U153").Value = "BED AND BREAKFAST"
$ws.Range("V153").Value = "Deluxe Room"
$ws.Range("W153").Value = "Single"
$ws.Range("X153").Value = "Refundable"
$ws.Range("Y153").Value = "Quiet Room"
$ws.Range("Z153").Value = "No Change"
$ws.Range("AA153").Formula = "=AD153-(AD153*25%)"
$ws.Range("AB153").Value = 0
$ws.Range("AC153").Value = 247.4
$ws.Range("AD153").Value = 2900
$ws.Range("AE153").Value = 2179.8000000000002
$ws.Range("AF153").Formula = "=AG153"
$ws.Range("AG153").Formula = "=AD153-AA153"
$ws.Range("AH153").Value = "Y"
$ws.Range("AI153").Formula = "=IF(O153<=3,""ShortStay"",""LongStay"")"
$ws.Range("AJ153").Value = "Leisure - Single"
$ws.Range("AK153").Value = 5
$ws.Range("AL153").Value = "Canadian"
$ws.Range("AM153").Value = 20
$ws.Range("AN153").Value = "Vancouver"
$ws.Range("AO153").Value = "Canada"
$ws.Range("AP153").Value = "Mobile"

# ---------------------------------------------------------------------------
# 5. Leave the selection where the workbook was last saved.
# ---------------------------------------------------------------------------
$ws.Range("AC147").Select()
